# Refresh the cryptos list (prices + 1h volume%) like the GitHub Actions
# job does on every run; also corrects rows 9/10 which had swapped
# Dogecoin / WrappedliquidstakedEther2.0 entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values like "0.9990" or
# "29.442.14" are not silently reinterpreted/rounded as numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.442.14'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.850.74'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D5').Value = '240.99'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').Value = '0.6328'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').Value = '1.0000'
$ws.Range('D8').Value = '3.901.91'
$ws.Range('E8').Value = '  +105.17%  '
$ws.Range('B9').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C9').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D9').Value = '4.140.95'
$ws.Range('E9').Value = '  +93.96%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.07566'
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('D11').Value = '0.2969'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '24.67'
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').Value = '0.07729'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '4.998'
$ws.Range('D15').Value = '0.6861'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '83.03'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '0.000009936'
$ws.Range('E17').Value = '  +4.16%  '
$ws.Range('D18').Value = '6.205'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '29.440.96'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').Value = '232.05'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').Value = '12.49'
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '7.605'
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '155.82'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('D26').Value = '0.1387'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '8.412'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = '17.69'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = '4.150.47'
$ws.Range('E29').Value = '  +102.62%  '
$ws.Range('D30').Value = '1.468'
$ws.Range('D31').Value = '0.05808'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').Value = '4.136'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('D34').Value = '4.024'
$ws.Range('D35').Value = '1.857'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').Value = '0.7174'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').Value = '2.596'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').Value = '1.251.86'
$ws.Range('E39').Value = '  +4.11%  '
$ws.Range('D40').Value = '2.804'
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('D41').Value = '0.01806'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').Value = '0.9035'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '101.58'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').Value = '67.12'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').Value = '7.206'
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('D48').Value = '9.150'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').Value = '1.690'
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('E51').Value = '  +0.06%  '
